$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values in row 1 (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON values in row 2 (B2:E2)
$ws.Range("B2").Value = 0.65556633424932542
$ws.Range("C2").Value = 1.4033779609428025
$ws.Range("D2").Value = 0.96313990293282503
$ws.Range("E2").Value = 1.3781050573816538

# Update STR values in row 3 (B3:E3)
$ws.Range("B3").Value = 1.544713802281203
$ws.Range("C3").Value = 0.73350108135927172
$ws.Range("D3").Value = 0.92422304724005855
$ws.Range("E3").Value = 1.0936601550003962

# Update the selection to match the committed state
$ws.Range("B1:E3").Select()
